$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph that follows the title (Heading1) paragraph.
$metaPara = $d.Paragraphs.Item(2)
[void]$metaPara.Range.Delete()

# 2. Replace the final paragraph (the italic image-prompt paragraph) with two paragraphs:
#    a new bold "Play Genius Free Slot Game by Cristaltec - Review" paragraph followed by
#    the (now updated) italic meta-description paragraph.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$xml = @"
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Genius Free Slot Game by Cristaltec - Review</w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Discover Genius, a top-tier online slot game by Cristaltec. Enjoy stunning graphics, versatile gameplay modes, and various ways to win for free.</w:t></w:r></w:p>
"@
[void]$lastPara.Range.InsertXML($xml)
